$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: id_end_use = 5, name = "ventilation"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "ventilation"

# Resize the table to include the new row
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:B6"))

# Set column A width to match the bestFit width applied by Excel
$ws.Columns.Item(1).ColumnWidth = 11.6666666666667

# Adjust view: zoom and selection
$ws.Application.ActiveWindow.Zoom = 261
$ws.Range("B7").Select()
